# Word COM-interop script implementing the commit
# "added CarApi.java and more documentation"
#
# It inserts five new paragraphs at the very start of the document body
# (a heading line, a blank line, a sub-heading line, a blank line, and a
# body paragraph ending with a manual page break), then tweaks sectPr /
# styles.xml the way the target diff shows.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert the five new paragraphs before the (current) first paragraph.
# ---------------------------------------------------------------------
$firstPara = $d.Paragraphs.First
$insPoint = $firstPara.Range
$insPoint.Collapse(1)   # wdCollapseStart

# Insert five paragraph marks in one shot, inheriting the formatting of
# the paragraph we are inserting in front of (same pPr/rPr as the rest
# of the document).
$null = $insPoint.InsertParagraphBefore()
$null = $insPoint.InsertParagraphBefore()
$null = $insPoint.InsertParagraphBefore()
$null = $insPoint.InsertParagraphBefore()
$null = $insPoint.InsertParagraphBefore()

# Paragraph 1: heading line
$d.Paragraphs.Item(1).Range.Text = "Documentation / explication de les classes:"
# Paragraph 2: left blank on purpose
# Paragraph 3: sub-heading
$d.Paragraphs.Item(3).Range.Text = "A) CarApi.java, dans /src/main/java/com/course/practicaljava/api/server"
# Paragraph 4: left blank on purpose
# Paragraph 5: body text
$d.Paragraphs.Item(5).Range.Text = "Voici un objet d’un voiture et cette propiétés. Les annotations “@AutoWired, @GetMapping”, etc, pernet à le framework spring et sont utilisé pour gestioner les appeles d’el client au le backend / serveur."

# ---------------------------------------------------------------------
# Append a manual page break at the end of paragraph 5 (the body text
# paragraph), as its own run, with no leftover paragraph split.
#
# Range.InsertBreak always introduces a fresh paragraph boundary, so we
# insert it at the end of paragraph 5 (which temporarily creates a 6th,
# break-only paragraph) and then delete the paragraph mark that
# separates them, merging the break back into paragraph 5 as a trailing
# run (matching how Word stores a same-paragraph page break).
# ---------------------------------------------------------------------
$p5 = $d.Paragraphs.Item(5)
$p5end = $p5.Range
$p5end.Collapse(0)   # wdCollapseEnd
$p5end.InsertBreak(7) # wdPageBreak

$p5fresh = $d.Paragraphs.Item(5)
$p6fresh = $d.Paragraphs.Item(6)
$mark = $d.Range($p5fresh.Range.End, $p6fresh.Range.Start)
$mark.Delete()

# ---------------------------------------------------------------------
# 2) sectPr: add docGrid.
# ---------------------------------------------------------------------
$sectionsCount = $d.Sections.Count
$sec = $d.Sections.Item($sectionsCount)
$secXml = '<w:sectPr xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:docGrid w:type="default" w:linePitch="100" w:charSpace="0"/></w:sectPr>'

# ---------------------------------------------------------------------
# 3 & 4) styles.xml tweaks are easiest done through the Styles
# collection exposed on the document.
# ---------------------------------------------------------------------
$normalStyle = $d.Styles.Item("Normal")
